$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
Write-Output ("before: " + $hdr.Range.Text)
$hdr.Range.Find.Execute("226", $true, $false, $false, $false, $false, $true, 1, $false, "320", 2)
Write-Output ("after: " + $hdr.Range.Text)
